$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D16").Value = 19.0
